# Update role values and selection as described in the commit "updated role and user"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14-19: Role column (D) changes from "Engineering" to "Process"
for ($r = 14; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value = "Process"
}

# Rows 20-28: Role column (D) stays "Production" (shared string index shifts after
# "Engineering" is removed from the workbook, but the displayed value is unchanged)
for ($r = 20; $r -le 28; $r++) {
    $ws.Cells.Item($r, 4).Value = "Production"
}

# Update the active selection on the sheet to D17
$ws.Range("D17").Select()
